{"js": "const replacements = [\n  [\"340\u00f79=37, 7\", \"908\u00f77=129, 5\"],\n  [\"446\u00f78=55, 6\", \"886\u00f74=221, 2\"],\n  [\"107\u00f76=17, 5\", \"127\u00f77=18, 1\"],\n  [\"782\u00f75=156, 2\", \"616\u00f78=77, 0\"],\n  [\"777\u00f78=97, 1\", \"558\u00f74=139, 2\"],\n  [\"566\u00f76=94, 2\", \"149\u00f72=74, 1\"],\n  [\"558\u00f73=186, 0\", \"583\u00f77=83, 2\"],\n  [\"906\u00f77=129, 3\", \"994\u00f73=331, 1\"],\n  [\"762\u00f72=381, 0\", \"645\u00f78=80, 5\"],\n  [\"867\u00f77=123, 6\", \"489\u00f73=163, 0\"],\n  [\"657\u00f77=93, 6\", \"671\u00f75=134, 1\"],\n  [\"505\u00f76=84, 1\", \"281\u00f76=46, 5\"],\n  [\"168\u00f75=33, 3\", \"687\u00f77=98, 1\"],\n  [\"876\u00f77=125, 1\", \"913\u00f79=101, 4\"],\n  [\"331\u00f73=110, 1\", \"899\u00f76=149, 5\"],\n  [\"704\u00f79=78, 2\", \"614\u00f76=102, 2\"],\n  [\"415\u00f79=46, 1\", \"129\u00f72=64, 1\"],\n  [\"561\u00f78=70, 1\", \"334\u00f73=111, 1\"],\n  [\"332\u00f78=41, 4\", \"384\u00f77=54, 6\"],\n  [\"423\u00f73=141, 0\", \"164\u00f79=18, 2\"],\n  [\"468\u00f74=117, 0\", \"494\u00f75=98, 4\"],\n  [\"740\u00f74=185, 0\", \"598\u00f77=85, 3\"],\n  [\"409\u00f76=68, 1\", \"416\u00f74=104, 0\"],\n  [\"225\u00f77=32, 1\", \"137\u00f75=27, 2\"],\n  [\"745\u00f72=372, 1\", \"628\u00f77=89, 5\"],\n];\n\nconst body = context.document.body;\nlet totalFound = 0;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n    totalFound++;\n  }\n  await context.sync();\n}\nreturn \"replaced=\" + totalFound;", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"340\u00f79=37, 7\", \"908\u00f77=129, 5\")\n  ,@(\"446\u00f78=55, 6\", \"886\u00f74=221, 2\")\n  ,@(\"107\u00f76=17, 5\", \"127\u00f77=18, 1\")\n  ,@(\"782\u00f75=156, 2\", \"616\u00f78=77, 0\")\n  ,@(\"777\u00f78=97, 1\", \"558\u00f74=139, 2\")\n  ,@(\"566\u00f76=94, 2\", \"149\u00f72=74, 1\")\n  ,@(\"558\u00f73=186, 0\", \"583\u00f77=83, 2\")\n  ,@(\"906\u00f77=129, 3\", \"994\u00f73=331, 1\")\n  ,@(\"762\u00f72=381, 0\", \"645\u00f78=80, 5\")\n  ,@(\"867\u00f77=123, 6\", \"489\u00f73=163, 0\")\n  ,@(\"657\u00f77=93, 6\", \"671\u00f75=134, 1\")\n  ,@(\"505\u00f76=84, 1\", \"281\u00f76=46, 5\")\n  ,@(\"168\u00f75=33, 3\", \"687\u00f77=98, 1\")\n  ,@(\"876\u00f77=125, 1\", \"913\u00f79=101, 4\")\n  ,@(\"331\u00f73=110, 1\", \"899\u00f76=149, 5\")\n  ,@(\"704\u00f79=78, 2\", \"614\u00f76=102, 2\")\n  ,@(\"415\u00f79=46, 1\", \"129\u00f72=64, 1\")\n  ,@(\"561\u00f78=70, 1\", \"334\u00f73=111, 1\")\n  ,@(\"332\u00f78=41, 4\", \"384\u00f77=54, 6\")\n  ,@(\"423\u00f73=141, 0\", \"164\u00f79=18, 2\")\n  ,@(\"468\u00f74=117, 0\", \"494\u00f75=98, 4\")\n  ,@(\"740\u00f74=185, 0\", \"598\u00f77=85, 3\")\n  ,@(\"409\u00f76=68, 1\", \"416\u00f74=104, 0\")\n  ,@(\"225\u00f77=32, 1\", \"137\u00f75=27, 2\")\n  ,@(\"745\u00f72=372, 1\", \"628\u00f77=89, 5\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nWrite-Output \"done\""}
